# PreCondition2.xlsx edit
#
# The author repointed the stimulus-image references from a POSIX-style
# relative path ("PreCondition/CS+2.BMP", "PreCondition/CS-2.BMP") to a
# Windows-style one ("PreCondition\CS+2.BMP", "PreCondition\CS-2.BMP"),
# matching the new absolute folder recorded elsewhere in the workbook
# (C:\Users\yhuang\Desktop\...\Test3\). Column A of Sheet1 holds these
# strings (shared across rows 2-6 and rows 7-11), so every cell that
# carries one of the two old values gets fixed up to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$old1 = "PreCondition/CS+2.BMP"
$new1 = "PreCondition\CS+2.BMP"
$old2 = "PreCondition/CS-2.BMP"
$new2 = "PreCondition\CS-2.BMP"

# Primary approach: a plain text find & replace across the used cells,
# same as a user doing Ctrl+H in Excel.
$ws.Cells.Replace($old1, $new1)
$ws.Cells.Replace($old2, $new2)

# Belt-and-braces fallback: walk the used range and fix up any cell that
# still carries the old value verbatim (covers cases where Replace()
# didn't touch a particular cell for any reason).
$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -eq $old1) {
            $cell.Value2 = $new1
        } elseif ($val -eq $old2) {
            $cell.Value2 = $new2
        }
    }
}
